$d = $word.ActiveDocument

# Locate the paragraph "ตาราง … Activity Diagram" (the 3rd paragraph in the
# document) and rework it into "ตารางที่ 1 Activity Diagram" while keeping
# the existing runs ("ตาราง" / "Activity Diagram") untouched and only
# rewriting the middle two runs (the single space, and the "… " run).

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("ตาราง")) {
        $target = $p
        break
    }
}

$pStart = $target.Range.Start

# Character offsets, relative to the paragraph start, of the run boundaries
# in the original text "ตาราง … Activity Diagram":
#   "ตาราง"            -> [0, 5)   (untouched)
#   " "                -> [5, 6)   (becomes "ที่ ")
#   "… "               -> [6, 8)   (becomes "1" then a new " " run)
#   "Activity Diagram" -> [8, …)   (untouched)

# --- Run 2: single space -> "ที่ " ------------------------------------
$spaceRng = $d.Range($pStart + 5, $pStart + 6)
# Nudge the formatting away from the neighbouring run first so the engine
# keeps this as its own run instead of silently merging it into "ตาราง"
# when the text is rewritten, then restore the original formatting.
$spaceRng.Font.Bold = $true
$spaceRng.Text = "ที่ "
$spaceRng2 = $d.Range($spaceRng.Start, $spaceRng.End)
$spaceRng2.Font.Bold = $false

# --- Run 3: "… " -> "1" -------------------------------------------------
$ellStart = $spaceRng2.End
$ellRng = $d.Range($ellStart, $ellStart + 2)
$ellRng.Font.Bold = $true
$ellRng.Text = "1"
$ellRng2 = $d.Range($ellRng.Start, $ellRng.End)
$ellRng2.Font.Bold = $false

# --- New run: " " inserted right after "1" -------------------------------
$insPoint = $d.Range($ellRng2.End, $ellRng2.End)
$insPoint.InsertAfter(" ")
$newSpace = $d.Range($ellRng2.End, $ellRng2.End + 1)
$newSpace.Font.Bold = $true
$newSpace2 = $d.Range($newSpace.Start, $newSpace.End)
$newSpace2.Font.Bold = $false
